$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 42602.016446759262

$ws.Range("B5").Value = "Bag"
$ws.Range("C5").Value = 4270
$ws.Range("D5").Value = 5159
$ws.Range("E5").Value = 602
$ws.Range("F5").Value = 105
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 71
$ws.Range("I5").Value = 27
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 0
